$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value without letting Excel auto-coerce numeric-looking
# strings (e.g. "23.00", "510.35") into real numbers, which would silently
# drop meaningful trailing zeros / change the stored type. The leading
# apostrophe forces Excel to treat the entry as literal text; resetting the
# style back to "Normal" afterwards clears the quote-prefix flag Excel
# stamps on the cell's style so no stray per-cell formatting is introduced.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue 2 4 "55.606.06"
Set-TextValue 2 5 "  -1.57%  "

# Row 3 - Ethereum
Set-TextValue 3 4 "2.283.84"
Set-TextValue 3 5 "  -1.31%  "

# Row 4 - TetherUSD
Set-TextValue 4 5 "  -0.05%  "

# Row 5 - BNB
Set-TextValue 5 4 "510.35"
Set-TextValue 5 5 "  -1.22%  "

# Row 6 - Solana
Set-TextValue 6 4 "128.68"
Set-TextValue 6 5 "  -4.69%  "

# Row 7 - USDC
Set-TextValue 7 5 "  +0.45%  "

# Row 8 - XRP
Set-TextValue 8 4 "0.524"
Set-TextValue 8 5 "  -2.26%  "

# Row 9 - LidoStakedEther
Set-TextValue 9 4 "2.282.59"
Set-TextValue 9 5 "  -2.16%  "

# Row 10 - Dogecoin
Set-TextValue 10 4 "0.0983"
Set-TextValue 10 5 "  -3.73%  "

# Row 11 - TRON
Set-TextValue 11 5 "  -0.22%  "

# Row 12 - Toncoin
Set-TextValue 12 5 "  -3.12%  "

# Row 13 - Cardano
Set-TextValue 13 4 "0.330"
Set-TextValue 13 5 "  -2.95%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue 14 4 "2.690.94"
Set-TextValue 14 5 "  -1.35%  "

# Row 15 - Avalanche
Set-TextValue 15 4 "23.00"
Set-TextValue 15 5 "  -3.96%  "

# Row 16 - WrappedBTC
Set-TextValue 16 4 "55.595.09"
Set-TextValue 16 5 "  -1.82%  "

# Row 17 - ShibaInu
Set-TextValue 17 5 "  -3.14%  "

# Row 18 - WrappedEther
Set-TextValue 18 4 "2.283.25"
Set-TextValue 18 5 "  -1.77%  "

# Row 19 - Chainlink
Set-TextValue 19 4 "10.20"
Set-TextValue 19 5 "  -2.78%  "

# Row 20 - BitcoinCash
Set-TextValue 20 4 "323.39"
Set-TextValue 20 5 "  +0.14%  "

# Row 21 - Polkadot
Set-TextValue 21 4 "4.07"
Set-TextValue 21 5 "  -3.71%  "

# Row 22 - Uniswap
Set-TextValue 22 4 "6.58"
Set-TextValue 22 5 "  +0.36%  "

# Row 23 - Dai
Set-TextValue 23 4 "0.999"

# Row 24 - Litecoin
Set-TextValue 24 4 "60.25"
Set-TextValue 24 5 "  -0.82%  "

# Row 25 - Kaspa
Set-TextValue 25 5 "  -1.18%  "

# Row 26 - Binance-PegBSC-USD
Set-TextValue 26 4 "0.997"
Set-TextValue 26 5 "  +0.55%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue 27 4 "8.42"

# Row 28 - Fetch.AI
Set-TextValue 28 5 "  +0.04%  "

# Row 29 - Monero
Set-TextValue 29 4 "166.99"
Set-TextValue 29 5 "  +0.10%  "

# Row 30 - PancakeSwap
Set-TextValue 30 5 "  -2.03%  "

# Row 31 - PEPE
Set-TextValue 31 5 "  -5.13%  "

# Row 32 - Aptos
Set-TextValue 32 4 "5.98"
Set-TextValue 32 5 "  -3.69%  "

# Row 33 - USDe
Set-TextValue 33 5 "  +0.03%  "

# Row 34 - EthereumClassic
Set-TextValue 34 4 "18.02"
Set-TextValue 34 5 "  -2.07%  "

# Row 35 - FirstDigitalUSD
Set-TextValue 35 4 "0.999"
Set-TextValue 35 5 "  +0.58%  "

# Row 36 - ImmutableX
Set-TextValue 36 5 "  -3.45%  "

# Row 37 - SuiNetwork
Set-TextValue 37 4 "0.875"
Set-TextValue 37 5 "  -4.76%  "

# Row 38 - NEARProtocol
Set-TextValue 38 4 "3.83"
Set-TextValue 38 5 "  -4.38%  "

# Row 39 - OKB
Set-TextValue 39 4 "38.17"
Set-TextValue 39 5 "  +0.55%  "

# Row 40 - Stacks
Set-TextValue 40 5 "  -0.72%  "

# Row 41 - Aave
Set-TextValue 41 4 "145.54"
Set-TextValue 41 5 "  +3.87%  "

# Row 42 - PolygonEcosystemToken
Set-TextValue 42 4 "0.369"
Set-TextValue 42 5 "  -3.52%  "

# Row 43 - now Filecoin (was Bittensor)
Set-TextValue 43 2 "Filecoin"
Set-TextValue 43 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue 43 4 "3.50"
Set-TextValue 43 5 "  -2.74%  "

# Row 44 - now Bittensor (was Filecoin)
Set-TextValue 44 2 "Bittensor"
Set-TextValue 44 3 "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue 44 4 "277.25"
Set-TextValue 44 5 "  +0.46%  "

# Row 45 - RenderToken
Set-TextValue 45 4 "4.84"
Set-TextValue 45 5 "  -6.28%  "

# Row 46 - Stellar
Set-TextValue 46 5 "  -1.85%  "

# Row 47 - Hedera
Set-TextValue 47 4 "0.0489"
Set-TextValue 47 5 "  -3.51%  "

# Row 48 - Mantle
Set-TextValue 48 4 "0.545"
Set-TextValue 48 5 "  -2.60%  "

# Row 49 - InjectiveProtocol
Set-TextValue 49 4 "17.67"
Set-TextValue 49 5 "  -0.47%  "

# Row 50 - Polygon
Set-TextValue 50 4 "0.376"
Set-TextValue 50 5 "  -1.34%  "

# Row 51 - VeChain
Set-TextValue 51 4 "0.0210"
Set-TextValue 51 5 "  -3.49%  "
